# Correction in SA algorithm and 746 logs
# Updates the Fitness values (column C) for run_17, generations 0-66
# (worksheet rows 2-68) to reflect the corrected algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block below is a contiguous range of rows whose Fitness value
# (column C) changes to a single corrected value, as described by the diff.
$blocks = @(
    @{ Start = 2;  End = 19; Value = 7736 },
    @{ Start = 20; End = 31; Value = 7345 },
    @{ Start = 32; End = 36; Value = 7343 },
    @{ Start = 37; End = 56; Value = 7295 },
    @{ Start = 57; End = 68; Value = 7293 }
)

foreach ($block in $blocks) {
    $range = $ws.Range("C$($block.Start):C$($block.End)")
    $range.Value = $block.Value
}
